$p = $ppt.ActivePresentation

# The deck ships with two DrawingML theme parts: the Integral theme (used
# by the slide master / main presentation theme) and the default Office
# Theme (used only by the notes master). The edit swaps which theme part
# carries which palette - the slide master's theme becomes the plain
# "Office Theme" colours while the notes-only theme becomes "Integral".
#
# The PowerPoint object model only exposes the palette that's attached to
# the slide master (Master.Theme.ThemeColorScheme), so reproduce the swap
# by writing the target ("Office Theme") RGB values into that scheme.
# RGB values use the classic VBA RGB(r,g,b) = r + g*256 + b*65536 packing.

$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
